# Add a new worksheet "JHSD200_11.2.28_24062025" at the end of the workbook,
# matching the structure/content added in the target diff (new sheet4.xml +
# one new shared string "JHSD200_11.2.28").

$wb = $excel.ActiveWorkbook

$sheetCount = $wb.Worksheets.Count
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($sheetCount))
$ws.Name = "JHSD200_11.2.28_24062025"

# Header row
$ws.Range("A1").Value = "DeviceName"
$ws.Range("B1").Value = "AppName"
$ws.Range("C1").Value = "ContentName"
$ws.Range("D1").Value = "Audio-Type"
$ws.Range("E1").Value = "Vision-Type"
$ws.Range("F1").Value = "VideoResolution"

# Data rows
$ws.Range("A2").Value = "JHSD200_11.2.28"
$ws.Range("B2").Value = "jiohotstar_25.06.02.3"
$ws.Range("C2").Value = "subham"
$ws.Range("D2").Value = "PCM"
$ws.Range("E2").Value = "src_fmt = SDR"
$ws.Range("F2").Value = "1920 1080"

$ws.Range("A3").Value = "JHSD200_11.2.28"
$ws.Range("B3").Value = "jiohotstar_25.06.02.3"
$ws.Range("C3").Value = "Salaar"
$ws.Range("D3").Value = "PCM"
$ws.Range("E3").Value = "src_fmt = SDR"
$ws.Range("F3").Value = "1920 1080"

$ws.Range("A4").Value = "JHSD200_11.2.28"
$ws.Range("B4").Value = "jiohotstar_25.06.02.3"
$ws.Range("C4").Value = "subham"
$ws.Range("D4").Value = "PCM"
$ws.Range("E4").Value = "src_fmt = SDR"
$ws.Range("F4").Value = "1920 1080"

$ws.Range("A5").Value = "JHSD200_11.2.28"
$ws.Range("B5").Value = "jiohotstar_25.06.02.3"
$ws.Range("C5").Value = "Salaar"
$ws.Range("D5").Value = "PCM"
$ws.Range("E5").Value = "src_fmt = SDR"
$ws.Range("F5").Value = "1920 1080"

$ws.Range("A6").Value = "JHSD200_11.2.28"
$ws.Range("B6").Value = "jiohotstar_25.06.02.3"
$ws.Range("C6").Value = "Bhagavanth Kesari"
$ws.Range("D6").Value = "PCM"
$ws.Range("E6").Value = "src_fmt = SDR"
$ws.Range("F6").Value = "1920 1080"

$ws.Range("A7").Value = "JHSD200_11.2.28"
$ws.Range("B7").Value = "jiohotstar_25.06.02.3"
$ws.Range("C7").Value = "Captain America Brave new world"
$ws.Range("D7").Value = "PCM"
$ws.Range("E7").Value = "src_fmt = SDR"
$ws.Range("F7").Value = "1920 1080"

$ws.Range("A8").Value = "JHSD200_11.2.28"
$ws.Range("B8").Value = "jiohotstar_25.06.02.3"
$ws.Range("C8").Value = "Anupama"
$ws.Range("D8").Value = "PCM"
$ws.Range("E8").Value = "src_fmt = SDR"
$ws.Range("F8").Value = "1920 1080"

# Column widths (best-fit-like, closest achievable to the authored widths)
$ws.Columns.Item(1).ColumnWidth = 14.0
$ws.Columns.Item(2).ColumnWidth = 17.333333333333336
$ws.Columns.Item(3).ColumnWidth = 27.333333333333336
$ws.Columns.Item(4).ColumnWidth = 9.333333333333332
$ws.Columns.Item(5).ColumnWidth = 11.333333333333332
$ws.Columns.Item(6).ColumnWidth = 13.333333333333332
